$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.811.23'
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '1.907.33'
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '

$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5225'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3784'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07242'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9062'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.74%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.937.69'
$ws.Range("E12").Value = '  +2.36%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07660'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.456'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008715'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.01%  '

$ws.Range("D19").Value = '27.840.44'
$ws.Range("E19").Value = '  +0.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.157'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.80%  '

$ws.Range("D22").Value = '2.153.16'
$ws.Range("E22").Value = '  +0.73%  '

$ws.Range("E23").Value = '  +1.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.635'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("E26").Value = '  -0.94%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.168'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.93'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.858'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09059'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.73%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.863'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.90%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.182'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.233'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7820'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02097'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.609'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.077'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5607'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.42%  '

$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05277'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.718'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '115.59'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.576'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.73%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4818'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.63%  '

$ws.Range("E47").Value = '  -1.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9997'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.622'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.06%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '66.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05991'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.01%  '
